# Generate Report for Archive
#
# The localization report workbook is regenerated and the four tracked
# files get re-ordered: the file that used to be last (c65f1799-...) is
# now reported first, and the other three (d424f980-..., f645ef34-...,
# 564bd0cf-...) each shift down one row. This touches the "Overview"
# sheet plus the per-language "zh-cn" / "de-de" detail sheets, and the
# hyperlinks that go with the file-name columns on each sheet.

$wb = $excel.ActiveWorkbook

function Read-RowValues($sheet, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $sheet.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Write-RowValues($sheet, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $sheet.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# Rows 2..5 hold one tracked file each; only rotate the columns whose
# content actually differs row-to-row (skip columns that are identical
# for every file, e.g. the file extension or the always-blank columns,
# so those cells are left exactly as they were).
function Rotate-Rows($sheet, $cols) {
    $r2 = Read-RowValues $sheet 2 $cols
    $r3 = Read-RowValues $sheet 3 $cols
    $r4 = Read-RowValues $sheet 4 $cols
    $r5 = Read-RowValues $sheet 5 $cols

    # new row2 = old row5, new row3 = old row2, new row4 = old row3, new row5 = old row4
    Write-RowValues $sheet 2 $cols $r5
    Write-RowValues $sheet 3 $cols $r2
    Write-RowValues $sheet 4 $cols $r3
    Write-RowValues $sheet 5 $cols $r4
}

function Rotate-Hyperlinks($sheet, $colLetter, $targets, $displays) {
    # $targets / $displays are ordered for rows 2..5 (post-rotation).
    $sheet.Hyperlinks.Delete()
    for ($i = 0; $i -lt 4; $i++) {
        $row = $i + 2
        $rng = $sheet.Range("$colLetter$row")
        $sheet.Hyperlinks.Add($rng, $targets[$i], "", "", $displays[$i])
    }
}

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Rotate-Rows $wsOverview @(1, 2, 5, 6, 7)

Rotate-Hyperlinks $wsOverview "B" @(
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77e840f978b528497761beeba5bef8baa7dc97d2/e2e/c65f1799-f719-4265-9023-c59b2653d4ed.md",
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9e5c5958cb6f59d9e2126b33604ce1348afce3f/e2e/d424f980-a6ca-4db1-844e-fc19a687f1a3.md",
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9e5c5958cb6f59d9e2126b33604ce1348afce3f/e2e/f645ef34-40ff-4d00-b57b-8f62574e529b.md",
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c11654a53fb7672a73c4fa7d8608ecc37a54e08b/e2e/564bd0cf-42e9-4340-a0e9-fb94fd5c91e4.md"
) @(
    "e2e\c65f1799-f719-4265-9023-c59b2653d4ed.md",
    "e2e\d424f980-a6ca-4db1-844e-fc19a687f1a3.md",
    "e2e\f645ef34-40ff-4d00-b57b-8f62574e529b.md",
    "e2e\564bd0cf-42e9-4340-a0e9-fb94fd5c91e4.md"
)

# ---- zh-cn / de-de detail sheets --------------------------------------
foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)
    Rotate-Rows $ws @(1, 3, 7, 8)

    Rotate-Hyperlinks $ws "A" @(
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77e840f978b528497761beeba5bef8baa7dc97d2/e2e/c65f1799-f719-4265-9023-c59b2653d4ed.md",
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9e5c5958cb6f59d9e2126b33604ce1348afce3f/e2e/d424f980-a6ca-4db1-844e-fc19a687f1a3.md",
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9e5c5958cb6f59d9e2126b33604ce1348afce3f/e2e/f645ef34-40ff-4d00-b57b-8f62574e529b.md",
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c11654a53fb7672a73c4fa7d8608ecc37a54e08b/e2e/564bd0cf-42e9-4340-a0e9-fb94fd5c91e4.md"
    ) @(
        "c65f1799-f719-4265-9023-c59b2653d4ed.md",
        "d424f980-a6ca-4db1-844e-fc19a687f1a3.md",
        "f645ef34-40ff-4d00-b57b-8f62574e529b.md",
        "564bd0cf-42e9-4340-a0e9-fb94fd5c91e4.md"
    )
}
